$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 109.3723396666667
$ws.Range("H2").Value = 328.117019
$ws.Range("I2").Value = 0.3006244632995563
$ws.Range("J2").Value = 0.3006244632995563
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 8.131233999999999
$ws.Range("N2").Value = 24.393702
$ws.Range("O2").Value = 0.02090995573015822
$ws.Range("P2").Value = 0.02090995573015823
$ws.Range("Q2").Value = 889.3320869571487
$ws.Range("R2").Value = 8003.988782614338
$ws.Range("S2").Value = 0.006286044218996299
$ws.Range("T2").Value = 0.0062860442189963

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 109.3723396666667
$ws.Range("H3").Value = 328.117019
$ws.Range("I3").Value = 0.3006244632995563
$ws.Range("J3").Value = 0.3006244632995563
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 243.3763986666667
$ws.Range("N3").Value = 730.1291960000001
$ws.Range("O3").Value = 0.625857000534647
$ws.Range("P3").Value = 0.6258570005346471
$ws.Range("Q3").Value = 26618.64614182075
$ws.Range("R3").Value = 239567.8152763868
$ws.Range("S3").Value = 0.1881479248879984
$ws.Range("T3").Value = 0.1881479248879984

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 109.3723396666667
$ws.Range("H4").Value = 328.117019
$ws.Range("I4").Value = 0.3006244632995563
$ws.Range("J4").Value = 0.3006244632995563
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 103.9426383333333
$ws.Range("N4").Value = 311.827915
$ws.Range("O4").Value = 0.2672947262403034
$ws.Range("P4").Value = 0.2672947262403035
$ws.Range("Q4").Value = 11368.44954564282
$ws.Range("R4").Value = 102316.0459107854
$ws.Range("S4").Value = 0.08035533361879306
$ws.Range("T4").Value = 0.08035533361879307

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 109.3723396666667
$ws.Range("H5").Value = 328.117019
$ws.Range("I5").Value = 0.3006244632995563
$ws.Range("J5").Value = 0.3006244632995563
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 33.41874933333333
$ws.Range("N5").Value = 100.256248
$ws.Range("O5").Value = 0.08593831749489127
$ws.Range("P5").Value = 0.08593831749489128
$ws.Range("Q5").Value = 3655.086803320523
$ws.Range("R5").Value = 32895.78122988471
$ws.Range("S5").Value = 0.02583516057376856
$ws.Range("T5").Value = 0.02583516057376856

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 170.751104
$ws.Range("H6").Value = 512.2533120000001
$ws.Range("I6").Value = 0.4693321835689973
$ws.Range("J6").Value = 0.4693321835689973
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 8.131233999999999
$ws.Range("N6").Value = 24.393702
$ws.Range("O6").Value = 0.02090995573015822
$ws.Range("P6").Value = 0.02090995573015823
$ws.Range("Q6").Value = 1388.417182382336
$ws.Range("R6").Value = 12495.75464144102
$ws.Range("S6").Value = 0.009813715181166227
$ws.Range("T6").Value = 0.009813715181166228

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 170.751104
$ws.Range("H7").Value = 512.2533120000001
$ws.Range("I7").Value = 0.4693321835689973
$ws.Range("J7").Value = 0.4693321835689973
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 243.3763986666667
$ws.Range("N7").Value = 730.1291960000001
$ws.Range("O7").Value = 0.625857000534647
$ws.Range("P7").Value = 0.6258570005346471
$ws.Range("Q7").Value = 41556.78875987747
$ws.Range("R7").Value = 374011.0988388972
$ws.Range("S7").Value = 0.293734832662869
$ws.Range("T7").Value = 0.293734832662869

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 170.751104
$ws.Range("H8").Value = 512.2533120000001
$ws.Range("I8").Value = 0.4693321835689973
$ws.Range("J8").Value = 0.4693321835689973
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 103.9426383333333
$ws.Range("N8").Value = 311.827915
$ws.Range("O8").Value = 0.2672947262403034
$ws.Range("P8").Value = 0.2672947262403035
$ws.Range("Q8").Value = 17748.32024808939
$ws.Range("R8").Value = 159734.8822328045
$ws.Range("S8").Value = 0.125450017522839
$ws.Range("T8").Value = 0.125450017522839

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 170.751104
$ws.Range("H9").Value = 512.2533120000001
$ws.Range("I9").Value = 0.4693321835689973
$ws.Range("J9").Value = 0.4693321835689973
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 33.41874933333333
$ws.Range("N9").Value = 100.256248
$ws.Range("O9").Value = 0.08593831749489127
$ws.Range("P9").Value = 0.08593831749489128
$ws.Range("Q9").Value = 5706.288342965931
$ws.Range("R9").Value = 51356.59508669338
$ws.Range("S9").Value = 0.04033361820212308
$ws.Range("T9").Value = 0.04033361820212309

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 68.96861966666667
$ws.Range("H10").Value = 206.905859
$ws.Range("I10").Value = 0.1895694499632422
$ws.Range("J10").Value = 0.1895694499632422
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 8.131233999999999
$ws.Range("N10").Value = 24.393702
$ws.Range("O10").Value = 0.02090995573015822
$ws.Range("P10").Value = 0.02090995573015823
$ws.Range("Q10").Value = 560.7999851666686
$ws.Range("R10").Value = 5047.199866500018
$ws.Range("S10").Value = 0.003963888806521838
$ws.Range("T10").Value = 0.003963888806521839

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 68.96861966666667
$ws.Range("H11").Value = 206.905859
$ws.Range("I11").Value = 0.1895694499632422
$ws.Range("J11").Value = 0.1895694499632422
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 243.3763986666667
$ws.Range("N11").Value = 730.1291960000001
$ws.Range("O11").Value = 0.625857000534647
$ws.Range("P11").Value = 0.6258570005346471
$ws.Range("Q11").Value = 16785.33427548438
$ws.Range("R11").Value = 151068.0084793594
$ws.Range("S11").Value = 0.1186433673469976
$ws.Range("T11").Value = 0.1186433673469976

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 68.96861966666667
$ws.Range("H12").Value = 206.905859
$ws.Range("I12").Value = 0.1895694499632422
$ws.Range("J12").Value = 0.1895694499632422
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 103.9426383333333
$ws.Range("N12").Value = 311.827915
$ws.Range("O12").Value = 0.2672947262403034
$ws.Range("P12").Value = 0.2672947262403035
$ws.Range("Q12").Value = 7168.780290361554
$ws.Range("R12").Value = 64519.022613254
$ws.Range("S12").Value = 0.05067091423144971
$ws.Range("T12").Value = 0.05067091423144972

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 68.96861966666667
$ws.Range("H13").Value = 206.905859
$ws.Range("I13").Value = 0.1895694499632422
$ws.Range("J13").Value = 0.1895694499632422
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 33.41874933333333
$ws.Range("N13").Value = 100.256248
$ws.Range("O13").Value = 0.08593831749489127
$ws.Range("P13").Value = 0.08593831749489128
$ws.Range("Q13").Value = 2304.845012506337
$ws.Range("R13").Value = 20743.60511255703
$ws.Range("S13").Value = 0.01629127957827301
$ws.Range("T13").Value = 0.01629127957827301

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 14.72510066666667
$ws.Range("H14").Value = 44.175302
$ws.Range("I14").Value = 0.0404739031682042
$ws.Range("J14").Value = 0.04047390316820419
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 8.131233999999999
$ws.Range("N14").Value = 24.393702
$ws.Range("O14").Value = 0.02090995573015822
$ws.Range("P14").Value = 0.02090995573015823
$ws.Range("Q14").Value = 119.7332391942227
$ws.Range("R14").Value = 1077.599152748004
$ws.Range("S14").Value = 0.0008463075234738604
$ws.Range("T14").Value = 0.0008463075234738604

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 14.72510066666667
$ws.Range("H15").Value = 44.175302
$ws.Range("I15").Value = 0.0404739031682042
$ws.Range("J15").Value = 0.04047390316820419
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 243.3763986666667
$ws.Range("N15").Value = 730.1291960000001
$ws.Range("O15").Value = 0.625857000534647
$ws.Range("P15").Value = 0.6258570005346471
$ws.Range("Q15").Value = 3583.741970257467
$ws.Range("R15").Value = 32253.6777323172
$ws.Range("S15").Value = 0.02533087563678203
$ws.Range("T15").Value = 0.02533087563678203

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 14.72510066666667
$ws.Range("H16").Value = 44.175302
$ws.Range("I16").Value = 0.0404739031682042
$ws.Range("J16").Value = 0.04047390316820419
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 103.9426383333333
$ws.Range("N16").Value = 311.827915
$ws.Range("O16").Value = 0.2672947262403034
$ws.Range("P16").Value = 0.2672947262403035
$ws.Range("Q16").Value = 1530.565813017259
$ws.Range("R16").Value = 13775.09231715533
$ws.Range("S16").Value = 0.01081846086722169
$ws.Range("T16").Value = 0.01081846086722169

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 14.72510066666667
$ws.Range("H17").Value = 44.175302
$ws.Range("I17").Value = 0.0404739031682042
$ws.Range("J17").Value = 0.04047390316820419
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 33.41874933333333
$ws.Range("N17").Value = 100.256248
$ws.Range("O17").Value = 0.08593831749489127
$ws.Range("P17").Value = 0.08593831749489128
$ws.Range("Q17").Value = 492.0944480874329
$ws.Range("R17").Value = 4428.850032786896
$ws.Range("S17").Value = 0.003478259140726618
$ws.Range("T17").Value = 0.003478259140726618
